$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Objetivos:" value (row 10) ---
$ws.Range("B10:C10").Value = "Propiciar uma integração entre os elementos de estruturação da cidade, das variáveis ambientais e da malha urbana."

# --- Insert a new row at 13 so that the "Docentes responsáveis:" value lives on its own
#     row, pushing "Programa resumido:" and everything below it down by one ---
$ws.Rows.Item(13).Insert()

# Bring over the correct formatting (wrap-text body styles) for the new B13:C13 cells
$ws.Range("B14:C14").Copy()
$ws.Range("B13:C13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# The inserted row copied the bold label style into A13 - the source data has no label there
$ws.Cells.Item(13, 1).Clear()

# Fill in the teacher name that belongs with "Docentes responsáveis:"
$ws.Range("B13:C13").Value = "5840942 - Marco Aurélio Kondracki de Alcântara"

# --- Fill in the rest of the values that were shifted down with the insert ---

# "Programa resumido:" (row 14 now)
$ws.Range("B14:C14").Value = "Variável Ecológicano Ambiente Urbano; Enfoque Encômico e Impactos Ambientais."

# "Programa:" (row 16 now)
$ws.Range("B16:C16").Value = "Elementos para estruturação ambiental da cidade. Variável ecológica no ambiente das atividades urbanas. A questão ambiental no urbanismo. A questão ambiental sob o enfoque econômico. Noções de higiene e saúde ambiental. A urbanização e os impactos ocasionados, principal enfoque da drenagem urbana."

# "Método:" (row 19 now)
$ws.Range("B19:C19").Value = "Aula expositiva e exercícios dirigidos."

# "Critério:" (row 20 now)
$ws.Range("B20:C20").Value = "Média ponderada de exercícios e provas."

# "Norma de recuperação:" (row 21 now)
$ws.Range("B21:C21").Value = "Prova única com nota igual ou superior a 5,0."

# "Bibliografia:" (row 22 now)
$ws.Range("B22:C22").Value = "valle, C.R. Qualidade ambiental: o desafio de ser competitivo protegendo o meio ambiente. Pioneira. 1995.`nDonaire, D.. Gestão ambiental na empresa. Atlas. 2a. edição. 1999.`nWinter, G.. Gestão e ambiente. Modelo prático de integração empresarial. Texto Editora, Lisboa. 1992.`nTucci, C.E., Porto, R.M., L.L. e Barros, M.T. org.. Drenagem Urbana. Ed. da Universidade e ABRH. 1995."
